# Refresh the crypto price list on Sheet1 (coin rows 2-51): update the
# "Price" (D) and "Volume(1h)" (E) columns with the latest scraped values.
#
# Several Price strings are plain-looking decimals (e.g. "213.49",
# "0.999", "0.530"). If assigned straight to .Value, Excel's COM layer
# auto-converts anything that parses as a number into a real number
# (e.g. "0.530" -> 0.53), which both changes the cell's stored type from
# text to numeric and silently drops the trailing zero. Since the sheet
# keeps every Price/Volume cell as literal text, for those ambiguous
# values we momentarily force Text number-format, assign the string,
# then restore the cell's original Style so the cell's formatting is
# left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.368.29"
$ws.Range("E2").Value = "  +3.03%  "

$ws.Range("D3").Value = "1.588.39"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("E4").Value = "  +1.08%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "213.49"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +1.30%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.492"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +0.25%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  +1.03%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.04"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +6.01%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.251"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +0.17%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0599"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +0.46%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0885"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("D12").Value = "1.814.27"
$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("D13").Value = "1.591.14"
$ws.Range("E13").Value = "  +1.80%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.530"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +1.92%  "

$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").Value = "28.347.07"
$ws.Range("E16").Value = "  +2.92%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.23"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +1.17%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "229.08"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").Value = "0.0₃0706"
$ws.Range("E19").Value = "  +0.16%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.48"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -0.53%  "

$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("E22").Value = "  -1.20%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.33"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("E24").Value = "  +0.48%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "151.87"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +1.38%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.18"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  -0.10%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.57"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("E29").Value = "  +1.04%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.14"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -0.06%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0472"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("E32").Value = "  +0.08%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.16"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("D34").Value = "1.396.16"
$ws.Range("E34").Value = "  -4.25%  "

$ws.Range("E35").Value = "  -1.50%  "

$ws.Range("E36").Value = "  -8.25%  "

$ws.Range("E37").Value = "  +1.62%  "

$ws.Range("E38").Value = "  -0.27%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.55"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +9.34%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.541"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -0.23%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.811"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -0.40%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  +1.08%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.89"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("E44").Value = "  -2.32%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.980"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +0.68%  "

$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("D47").Value = "1.723.46"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("E48").Value = "  +1.65%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "87.04"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("D50").Value = "0.0₆0104"
$ws.Range("E50").Value = "  +15.12%  "

$ws.Range("E51").Value = "  -0.78%  "
